# Update "想去人数" (F column) counts on both the "展览" sheet and the
# "全部类型" aggregate sheet, as captured in the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Map: row number -> new value, for the "展览" worksheet.
$zhanlanUpdates = @{
    2  = 165
    6  = 61
    9  = 2298
    10 = 115
    11 = 64
    12 = 154
    13 = 1402
    14 = 498
    15 = 32
    16 = 309
    17 = 215
    18 = 14
    23 = 3
    24 = 72
    25 = 29
    26 = 1432
    29 = 149
    30 = 180
    31 = 283
}

# Map: row number -> new value, for the "全部类型" worksheet.
$quanbuUpdates = @{
    2  = 165
    7  = 61
    10 = 2298
    11 = 115
    12 = 64
    13 = 154
    14 = 1402
    15 = 498
    16 = 32
    17 = 309
    18 = 215
    19 = 14
    24 = 3
    25 = 72
    26 = 29
    27 = 1432
    30 = 149
    31 = 180
    32 = 283
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Range("F$row").Value = $zhanlanUpdates[$row]
}

$wsQuanbu = $wb.Worksheets.Item("全部类型")
foreach ($row in $quanbuUpdates.Keys) {
    $wsQuanbu.Range("F$row").Value = $quanbuUpdates[$row]
}
